$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$metadata = $wb.Worksheets.Item("Metadata")

# URL
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/claim-supporting-info-type"

# Version
$metadata.Range("B3").Value = "8.0.0"

# Date
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$metadata.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Claim Supporting" ---
$codes = $wb.Worksheets.Item("Include from Claim Supporting")

# System URI
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/claim-supporting-info-type"
